$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Paragraphs.Item(1).Range.Text = "2024-08-27 Tuesday"

# Update the division problems in the table (5 columns, data rows 1,5,9,13,17)
$t = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)

$values = @(
  @("59÷6=", "46÷9=", "75÷8=", "43÷4=", "93÷2="),
  @("71÷5=", "58÷5=", "96÷6=", "13÷3=", "95÷4="),
  @("33÷8=", "83÷5=", "16÷8=", "36÷7=", "59÷5="),
  @("43÷7=", "25÷9=", "87÷7=", "72÷5=", "92÷8="),
  @("32÷5=", "48÷6=", "23÷7=", "30÷4=", "84÷2=")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $row = $rows[$i]
  for ($col = 1; $col -le 5; $col++) {
    $t.Cell($row, $col).Range.Text = $values[$i][$col - 1]
  }
}

Write-Output "done"
